# Update the cryptocurrency price table on the active worksheet with the
# latest scraped values (coin prices, 1h volume %, and two re-ordered rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The Price and Volume(1h) columns store plain numeric-looking text (e.g.
# "1.004", "0.4440", "  +0.07%  "). Pre-format them as Text so Excel does not
# silently coerce the assigned strings into numbers (which would also drop
# significant trailing zeros).
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "28.336.67"
$ws.Range("E2").Value = "  +0.97%  "
$ws.Range("D3").Value = "1.802.87"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "326.75"
$ws.Range("E5").Value = "  -3.21%  "
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("D7").Value = "0.4440"
$ws.Range("E7").Value = "  +4.85%  "
$ws.Range("D8").Value = "0.3755"
$ws.Range("E8").Value = "  +6.11%  "
$ws.Range("D9").Value = "44.69"
$ws.Range("E9").Value = "  -1.69%  "
$ws.Range("D10").Value = "1.149"
$ws.Range("E10").Value = "  -0.96%  "
$ws.Range("D11").Value = "0.07510"
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").Value = "22.57"
$ws.Range("E12").Value = "  -1.31%  "
$ws.Range("D13").Value = "0.9988"
$ws.Range("E13").Value = "  -0.33%  "
$ws.Range("D14").Value = "7.671"
$ws.Range("E14").Value = "  +5.12%  "
$ws.Range("D15").Value = "6.301"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").Value = "1.801.93"
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("D17").Value = "0.00001091"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").Value = "0.06800"
$ws.Range("E18").Value = "  +1.59%  "
$ws.Range("D19").Value = "80.79"
$ws.Range("E19").Value = "  -2.09%  "
$ws.Range("D20").Value = "0.9996"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("D22").Value = "6.313"
$ws.Range("E22").Value = "  -1.09%  "
$ws.Range("D23").Value = "28.363.74"
$ws.Range("E23").Value = "  +0.97%  "
$ws.Range("D24").Value = "11.81"
$ws.Range("E24").Value = "  -0.67%  "
$ws.Range("D25").Value = "2.409"
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").Value = "20.53"
$ws.Range("E26").Value = "  -1.33%  "
$ws.Range("D27").Value = "153.52"
$ws.Range("E27").Value = "  -1.68%  "
$ws.Range("D28").Value = "2.360"
$ws.Range("E28").Value = "  -5.70%  "
$ws.Range("D29").Value = "2.009.79"
$ws.Range("E29").Value = "  -0.65%  "
$ws.Range("D30").Value = "132.87"
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").Value = "1.254"
$ws.Range("E31").Value = "  -4.86%  "
$ws.Range("D32").Value = "4.001"
$ws.Range("E32").Value = "  -2.04%  "
$ws.Range("D33").Value = "5.823"
$ws.Range("E33").Value = "  -3.09%  "
$ws.Range("D34").Value = "0.09320"
$ws.Range("E34").Value = "  +1.36%  "
$ws.Range("D35").Value = "0.2284"
$ws.Range("E35").Value = "  +5.53%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.06360"
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").Value = "12.10"
$ws.Range("E37").Value = "  -2.45%  "
$ws.Range("D38").Value = "0.02317"
$ws.Range("E38").Value = "  -1.54%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "0.6581"
$ws.Range("E39").Value = "  -1.50%  "
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").Value = "5.156"
$ws.Range("E40").Value = "  -1.83%  "
$ws.Range("D41").Value = "1.206"
$ws.Range("E41").Value = "  -1.29%  "
$ws.Range("D42").Value = "1.452"
$ws.Range("E42").Value = "  -4.16%  "
$ws.Range("D43").Value = "8.144"
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("D44").Value = "0.9989"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").Value = "13.92"
$ws.Range("E45").Value = "  -1.77%  "
$ws.Range("D46").Value = "0.6068"
$ws.Range("E46").Value = "  -1.77%  "
$ws.Range("D47").Value = "3.795"
$ws.Range("E47").Value = "  -2.13%  "
$ws.Range("D48").Value = "128.72"
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("D49").Value = "2.032"
$ws.Range("E49").Value = "  -1.39%  "
$ws.Range("D50").Value = "0.07089"
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("D51").Value = "1.155"
$ws.Range("E51").Value = "  -2.66%  "
